$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Header cell F1: value + formatting copied from E1 (style s="1")
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Timestamps for F2:F65 (no special style, same as original column data cells)
$timestamps = @(
    "2021-10-05 13:41:47.061360",
    "2021-10-05 13:41:47.061373",
    "2021-10-05 13:41:47.061377",
    "2021-10-05 13:41:47.061380",
    "2021-10-05 13:41:47.061384",
    "2021-10-05 13:41:47.061387",
    "2021-10-05 13:41:47.061390",
    "2021-10-05 13:41:47.061393",
    "2021-10-05 13:41:47.061396",
    "2021-10-05 13:41:47.061399",
    "2021-10-05 13:41:47.061402",
    "2021-10-05 13:41:47.061405",
    "2021-10-05 13:41:47.061408",
    "2021-10-05 13:41:47.061410",
    "2021-10-05 13:41:47.061413",
    "2021-10-05 13:41:47.061416",
    "2021-10-05 13:41:47.061419",
    "2021-10-05 13:41:47.061422",
    "2021-10-05 13:41:47.061425",
    "2021-10-05 13:41:47.061428",
    "2021-10-05 13:41:47.061431",
    "2021-10-05 13:41:47.061434",
    "2021-10-05 13:41:47.061437",
    "2021-10-05 13:41:47.061440",
    "2021-10-05 13:41:47.061443",
    "2021-10-05 13:41:47.061446",
    "2021-10-05 13:41:47.061449",
    "2021-10-05 13:41:47.061452",
    "2021-10-05 13:41:47.061455",
    "2021-10-05 13:41:47.061458",
    "2021-10-05 13:41:47.061461",
    "2021-10-05 13:41:47.061464",
    "2021-10-05 13:41:47.061467",
    "2021-10-05 13:41:47.061470",
    "2021-10-05 13:41:47.061473",
    "2021-10-05 13:41:47.061476",
    "2021-10-05 13:41:47.061479",
    "2021-10-05 13:41:47.061482",
    "2021-10-05 13:41:47.061485",
    "2021-10-05 13:41:47.061488",
    "2021-10-05 13:41:47.061491",
    "2021-10-05 13:41:47.061494",
    "2021-10-05 13:41:47.061497",
    "2021-10-05 13:41:47.061500",
    "2021-10-05 13:41:47.061503",
    "2021-10-05 13:41:47.061506",
    "2021-10-05 13:41:47.061509",
    "2021-10-05 13:41:47.061512",
    "2021-10-05 13:41:47.061515",
    "2021-10-05 13:41:47.061518",
    "2021-10-05 13:41:47.061521",
    "2021-10-05 13:41:47.061524",
    "2021-10-05 13:41:47.061527",
    "2021-10-05 13:41:47.061530",
    "2021-10-05 13:41:47.061533",
    "2021-10-05 13:41:47.061536",
    "2021-10-05 13:41:47.061539",
    "2021-10-05 13:41:47.061541",
    "2021-10-05 13:41:47.061544",
    "2021-10-05 13:41:47.061547",
    "2021-10-05 13:41:47.061550",
    "2021-10-05 13:41:47.061553",
    "2021-10-05 13:41:47.061556",
    "2021-10-05 13:41:47.061559"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

